# "new changes in ops (ordercreation & orderpage & order form)"
#
# Content updates to the "with_all_correctdata" sheet:
#  - Row 2 order ST18-002 -> ST18-001, received-date shifted a day earlier,
#    and the product name "Full Search" -> "Current Owner Search".
#  - Row 4 order's Emp ID / Assignee_QA (SIPL5316/SIPL5688) cleared and
#    re-assigned to Typist/Typist QC SIPL5317, status WIP -> Typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates -------------------------------------------------------
$ws.Range("A2").Value = 45436.041666666664
$ws.Range("B2").Value = "ST18-001"
$ws.Range("J2").Value = "Current Owner Search"

# --- Row 4 updates -------------------------------------------------------
# Clear the Emp ID-Order Assigned / Assignee_QA values (column C/D) ...
$ws.Range("C4:D4").ClearContents()

# ... and populate the Typist / Typist QC columns (E/F) instead, copying
# the number formatting/border style from a neighbouring formatted cell so
# the new cells keep the same cell style as the rest of the row.
$ws.Range("H4").Copy() | Out-Null
$ws.Range("E4:F4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("E4").Value = "SIPL5317"
$ws.Range("F4").Value = "SIPL5317"

# Status moves from WIP to Typing for this order.
$ws.Range("M4").Value = "Typing"

# --- Cosmetic updates ------------------------------------------------------
# Column widths were re-autofit by Excel after the content changes above.
$ws.Columns.Item(3).ColumnWidth = 19.666666666666668
$ws.Columns.Item(7).ColumnWidth = 6.833333333333333
$ws.Columns.Item(10).ColumnWidth = 18.333333333333332

# Selection moved to H8 in the saved file.
$ws.Range("H8").Select() | Out-Null
